# "Capitalizar autor y Uppercase Título"
# Mark the corresponding Feature rows as done ("OK") and update the
# current cell selection, mirroring the author's manual edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  -> "Crear mas de un autor..."        -> OK in column F
$ws.Range("F2").Value = "OK"

# Row 4  -> "Convertir en mayuscula el título" -> OK in column C
$ws.Range("C4").Value = "OK"

# Row 5  -> "Capitalizar campos de autor"      -> OK in column C
$ws.Range("C5").Value = "OK"

# Row 11 -> "cambiar campo nota..."            -> OK in column C
$ws.Range("C11").Value = "OK"

# Row 12 -> "Crear autor"                      -> OK in column C
$ws.Range("C12").Value = "OK"

# Best-effort: the saved workbook window was minimized by the author;
# try to reflect that in the window state (no-op if unsupported).
try {
    $wb.Windows.Item(1).WindowState = -4140
} catch {
}

# Leave the selection where the author left it after editing.
$ws.Range("F4").Select()
